# Updated symbol list on Fri Dec 16 21:23:26 UTC 2022 with GitHub Actions
#
# This script updates the "cryptos" worksheet with refreshed price data
# (column D), reorders two rows whose coin info got swapped (rows 42/43),
# and tweaks a couple of the "Worst in 24h" suffixes on column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must stay stored as TEXT
# (matching the source workbook, which uses inline/shared strings for
# the whole column). Setting NumberFormat to "@" (Text) before writing
# the value prevents Excel from auto-converting the string to a number.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Column D price refresh -------------------------------------------------
Set-TextValue "D2"  "243.51"
Set-TextValue "D3"  "23.18"
Set-TextValue "D4"  "5.750"
Set-TextValue "D5"  "0.05808"
Set-TextValue "D6"  "3.420"
Set-TextValue "D7"  "6.473"
Set-TextValue "D8"  "1.321"
Set-TextValue "D9"  "0.8008"
Set-TextValue "D10" "0.1464"
Set-TextValue "D11" "0.07707"
Set-TextValue "D12" "0.03239"
Set-TextValue "D13" "0.03009"
Set-TextValue "D14" "0.09245"
Set-TextValue "D15" "0.001661"
Set-TextValue "D16" "3.253"
Set-TextValue "D17" "0.04763"
Set-TextValue "D18" "0.0005998"
Set-TextValue "D20" "0.005393"
Set-TextValue "D23" "3.689"
Set-TextValue "D25" "0.3323"
Set-TextValue "D26" "0.1242"
Set-TextValue "D27" "0.001001"
Set-TextValue "D40" "0.04297"
Set-TextValue "D41" "0.007076"
Set-TextValue "D44" "0.009724"
Set-TextValue "D46" "0.00005626"
Set-TextValue "D48" "0.7863"
Set-TextValue "D49" "0.09934"
Set-TextValue "D50" "0.00002103"

# --- Column E text tweaks ----------------------------------------------------
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E45").Value = "44ACDXExchangeACXTWorstin24h"

# --- Rows 42 / 43: CEJI and BKEXToken swapped places ------------------------
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1058"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003447"
$ws.Range("E43").Value = "42CEJICEJI"
